$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1, J1 - copy header style from H1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A1").Select()

$iValues = @(9,8,8,8,7,8,6,8,6,4,5,8,7,8,1,8,6,8,8,8,4,7,4,8,9,6,6,6,7,4,7,8,6,9,9,5,9,7,7,8,8,7,8,7,8,7,8,6,9,9,8,8,9,9,9,8,8,9,9,3,4,6,6,5,4,4,4)
$jValues = @(9,8,8,8,7,8,6,8,6,5,7,8,7,9,2,8,7,8,8,8,6,8,5,8,9,6,6,6,7,4,8,8,7,9,9,6,9,7,7,8,8,7,8,7,8,7,8,6,9,9,8,8,9,9,9,9,8,9,9,3,4,6,6,5,4,4,4)

for ($r = 2; $r -le 68; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $iValues[$idx]
    $ws.Cells.Item($r, 10).Value = $jValues[$idx]
}
